# BigAsciiChars.xlsx — "WIP of conversion to class based font"
#
# Renames the worksheet and redraws the 5x5 glyph-grid (columns J:N, rows
# 2:6) that the sheet uses to derive a packed font bitmap value. The glyph
# cells L4/M4/N5/L6/M6/N6 are plain inputs; columns O and P recompute
# automatically via their existing shared formulas once the inputs change,
# which also ripples into the P7 grand-total formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tab was renamed from the generic "Sheet1" to "5x5" (describing the
# pixel grid used to build each character).
$ws.Name = "5x5"

# Redraw the glyph: clear L4/M4, and turn on N5, L6, M6, N6.
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 1

# Move the on-sheet selection to the column that was just edited.
$ws.Range("P2:P6").Select()
